$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing row 74 (01-01-2021) values for columns B..R
$row74 = @(11048, 8324, 7248, 981, 95, 5232, 3288, 1945, 1202, 743, -3690, -274, 970, -37, -1181, -26, 1455)
for ($i = 0; $i -lt $row74.Length; $i++) {
    $col = $i + 2  # Column B = 2
    $ws.Cells.Item(74, $col).Value = $row74[$i]
}

# Add new row 75 (01-04-2021)
$ws.Cells.Item(75, 1).Value = "'01-04-2021"
$ws.Cells.Item(75, 1).ClearFormats()
$row75 = @(2491, 1448, 466, 785, 196, -5356, -3088, -2268, 967, -3236, -2731, 4685, 1298, 78, 3232, 77, 4446)
for ($i = 0; $i -lt $row75.Length; $i++) {
    $col = $i + 2  # Column B = 2
    $ws.Cells.Item(75, $col).Value = $row75[$i]
}
